$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row - add new column F "Modelo", matching style of existing headers (A1:E1)
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F1").Value = "Modelo"

# Update numeric values in row 2
$ws.Range("B2").Value = 0.6645926274483194
$ws.Range("C2").Value = 0.9804791735888083
$ws.Range("D2").Value = 0.6051733866243391

# Add new model description cell, with embedded newline
$modelo = "Pipeline(steps=[('model'," + [char]10 + "                 AdaBoostRegressor(learning_rate=0.5, n_estimators=100))])"
$ws.Range("F2").Value = $modelo
